$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.250.24'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').Value = '3.567.31'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''609.85'
$ws.Range('E5').Value = '  +3.61%  '
$ws.Range('D6').Value = '''186.11'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').Value = '3.562.22'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('E10').Value = '  +7.56%  '
$ws.Range('D11').Value = '''0.644'
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('D12').Value = '''53.84'
$ws.Range('E12').Value = '  -1.33%  '
$ws.Range('E13').Value = '  +0.67%  '
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').Value = '4.133.29'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '70.324.05'
$ws.Range('E16').Value = '  +0.00%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.573.82'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').Value = '''12.69'
$ws.Range('E18').Value = '  +1.62%  '
$ws.Range('D19').Value = '''18.97'
$ws.Range('E19').Value = '  -2.60%  '
$ws.Range('D20').Value = '''580.55'
$ws.Range('E20').Value = '  +6.65%  '
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('D22').Value = '''0.994'
$ws.Range('E22').Value = '  -2.52%  '
$ws.Range('D23').Value = '''17.34'
$ws.Range('E23').Value = '  -3.61%  '
$ws.Range('D24').Value = '''4.72'
$ws.Range('E24').Value = '  +0.56%  '
$ws.Range('E25').Value = '  -0.77%  '
$ws.Range('D26').Value = '''94.89'
$ws.Range('E26').Value = '  -1.17%  '
$ws.Range('E27').Value = '  -1.56%  '
$ws.Range('D28').Value = '''10.97'
$ws.Range('E28').Value = '  -2.52%  '
$ws.Range('D29').Value = '''9.39'
$ws.Range('E29').Value = '  +2.72%  '
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('D31').Value = '''7.04'
$ws.Range('E31').Value = '  -4.23%  '
$ws.Range('D32').Value = '''12.25'
$ws.Range('E32').Value = '  -2.39%  '
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('D35').Value = '''3.69'
$ws.Range('E35').Value = '  +19.51%  '
$ws.Range('D36').Value = '''3.24'
$ws.Range('E36').Value = '  +0.51%  '
$ws.Range('D37').Value = '''531.96'
$ws.Range('E37').Value = '  -4.32%  '
$ws.Range('D38').Value = '''0.403'
$ws.Range('E38').Value = '  -3.13%  '
$ws.Range('E39').Value = '  +0.16%  '
$ws.Range('D40').Value = '''37.13'
$ws.Range('E40').Value = '  -3.84%  '
$ws.Range('D41').Value = '0.0₃0785'
$ws.Range('E41').Value = '  +2.22%  '
$ws.Range('D42').Value = '3.539.24'
$ws.Range('E42').Value = '  +5.17%  '
$ws.Range('D43').Value = '''3.53'
$ws.Range('E43').Value = '  +3.84%  '
$ws.Range('E44').Value = '  +0.89%  '
$ws.Range('E45').Value = '  +3.59%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '''3.46'
$ws.Range('E46').Value = '  -2.70%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D47').Value = '''2.93'
$ws.Range('E47').Value = '  -1.72%  '
$ws.Range('E48').Value = '  +2.92%  '
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('D51').Value = '''136.42'
$ws.Range('E51').Value = '  -0.73%  '
